$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1975
$ws.Range("I18").Value = 1950
$ws.Range("K18").Value = 1950
$ws.Range("M18").Value = -1666

$ws.Range("H20").Value = 1021
$ws.Range("I20").Value = 1021
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 1021
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -791
$ws.Range("N20").ClearContents()

$ws.Range("H35").Value = 1021
$ws.Range("I35").Value = 1021
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 1021
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -642
$ws.Range("N35").ClearContents()

$ws.Range("H40").Value = 2227.5
$ws.Range("J40").Value = 1875
$ws.Range("L40").Value = 1875
$ws.Range("N40").Value = -2225

$ws.Range("H51").Value = 6156.25
$ws.Range("I51").Value = 3000
$ws.Range("J51").Value = 6366.6665
$ws.Range("K51").Value = 3000
$ws.Range("L51").Value = 6366.6665
$ws.Range("M51").Value = -2516
$ws.Range("N51").Value = -7334.6665

$ws.Range("H74").Value = 4632.706
$ws.Range("J74").Value = 4686
$ws.Range("L74").Value = 4686
$ws.Range("N74").Value = -6558

$ws.Range("H77").Value = 4632.706
$ws.Range("J77").Value = 4686
$ws.Range("L77").Value = 23430
$ws.Range("N77").Value = -32790

$ws.Range("H127").Value = 20409932
$ws.Range("I127").Value = 76923544
$ws.Range("J127").Value = 2240.111
$ws.Range("K127").Value = 230770632
$ws.Range("L127").Value = 6720.333
$ws.Range("M127").Value = -230765672
$ws.Range("N127").Value = -16640.333

$ws.Range("H129").Value = 894.97675
$ws.Range("J129").Value = 915.9459000000001
$ws.Range("L129").Value = 2747.8377
$ws.Range("N129").Value = -12747.8377

$ws.Range("H138").Value = 2063.6316
$ws.Range("I138").Value = 1721.1538
$ws.Range("J138").Value = 2805.6667
$ws.Range("K138").Value = 5163.4614
$ws.Range("L138").Value = 8417.000100000001
$ws.Range("M138").Value = -23.46140000000014
$ws.Range("N138").Value = -18697.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1590.2222
$ws.Range("I45").Value = 1065.0714
$ws.Range("K45").Value = 1065.0714
$ws.Range("M45").Value = -688.0714

$ws.Range("H63").Value = 8715.200000000001
$ws.Range("J63").Value = 2199.5
$ws.Range("L63").Value = 2199.5
$ws.Range("N63").Value = -3571.5

$ws.Range("H66").Value = 8715.200000000001
$ws.Range("J66").Value = 2199.5
$ws.Range("L66").Value = 10997.5
$ws.Range("N66").Value = -17861.5

$ws.Range("H102").Value = 1205
$ws.Range("I102").Value = 1205
$ws.Range("K102").Value = 1205
$ws.Range("M102").Value = 417

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 515.8
$ws.Range("I64").Value = 519.75
$ws.Range("J64").Value = 500
$ws.Range("K64").Value = 519.75
$ws.Range("L64").Value = 500
$ws.Range("M64").Value = -294.75
$ws.Range("N64").Value = -950

$ws.Range("H67").Value = 515.8
$ws.Range("I67").Value = 519.75
$ws.Range("J67").Value = 500
$ws.Range("K67").Value = 519.75
$ws.Range("L67").Value = 500
$ws.Range("M67").Value = 260.25
$ws.Range("N67").Value = -2060

$ws.Range("H82").Value = 29000
$ws.Range("I82").Value = 20000
$ws.Range("J82").Value = 38000
$ws.Range("K82").Value = 20000
$ws.Range("L82").Value = 38000
$ws.Range("M82").Value = -19617
$ws.Range("N82").Value = -38766

$ws.Range("H85").Value = 29000
$ws.Range("I85").Value = 20000
$ws.Range("J85").Value = 38000
$ws.Range("K85").Value = 20000
$ws.Range("L85").Value = 38000
$ws.Range("M85").Value = -18674
$ws.Range("N85").Value = -40652

$ws.Range("H86").Value = 2226.4075
$ws.Range("I86").Value = 2175.95
$ws.Range("J86").Value = 2370.5715
$ws.Range("K86").Value = 2175.95
$ws.Range("L86").Value = 2370.5715
$ws.Range("M86").Value = -1052.95
$ws.Range("N86").Value = -4616.5715

$ws.Range("H89").Value = 2226.4075
$ws.Range("I89").Value = 2175.95
$ws.Range("J89").Value = 2370.5715
$ws.Range("K89").Value = 10879.75
$ws.Range("L89").Value = 11852.8575
$ws.Range("M89").Value = -5263.75
$ws.Range("N89").Value = -23084.8575

$ws.Range("H105").Value = 3169.8333
$ws.Range("I105").Value = 3169.8333
$ws.Range("K105").Value = 3169.8333
$ws.Range("M105").Value = -1422.8333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1731.9
$ws.Range("I16").Value = 1731.9
$ws.Range("K16").Value = 1731.9
$ws.Range("M16").Value = -1444.9

$ws.Range("H31").Value = 14858.308
$ws.Range("I31").Value = 2940.5
$ws.Range("J31").Value = 25073.572
$ws.Range("K31").Value = 2940.5
$ws.Range("L31").Value = 25073.572
$ws.Range("M31").Value = -2645.5
$ws.Range("N31").Value = -25663.572

$ws.Range("H34").Value = 14858.308
$ws.Range("I34").Value = 2940.5
$ws.Range("J34").Value = 25073.572
$ws.Range("K34").Value = 2940.5
$ws.Range("L34").Value = 25073.572
$ws.Range("M34").Value = -2738.5
$ws.Range("N34").Value = -25477.572

$ws.Range("H86").Value = 2741.8667
$ws.Range("I86").Value = 2763.077
$ws.Range("K86").Value = 2763.077
$ws.Range("M86").Value = -1640.077

$ws.Range("H89").Value = 2741.8667
$ws.Range("I89").Value = 2763.077
$ws.Range("K89").Value = 13815.385
$ws.Range("M89").Value = -8199.385000000002

$ws.Range("H113").Value = 1731.9
$ws.Range("I113").Value = 1731.9
$ws.Range("K113").Value = 1731.9
$ws.Range("M113").Value = 438.0999999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1051.0968
$ws.Range("I5").Value = 556
$ws.Range("J5").Value = 1515.25
$ws.Range("K5").Value = 1668
$ws.Range("L5").Value = 4545.75
$ws.Range("M5").Value = -1556
$ws.Range("N5").Value = -4769.75

$ws.Range("H18").Value = 733.9524
$ws.Range("I18").Value = 534.73334
$ws.Range("J18").Value = 1232
$ws.Range("K18").Value = 1604.20002
$ws.Range("L18").Value = 3696
$ws.Range("M18").Value = -1435.20002
$ws.Range("N18").Value = -4034

$ws.Range("H131").Value = 1101.1666
$ws.Range("I131").Value = 365
$ws.Range("J131").Value = 1223.8611
$ws.Range("K131").Value = 1095
$ws.Range("L131").Value = 3671.5833
$ws.Range("M131").Value = 3945
$ws.Range("N131").Value = -13751.5833

$ws.Range("H135").Value = 1051.0968
$ws.Range("I135").Value = 556
$ws.Range("J135").Value = 1515.25
$ws.Range("K135").Value = 5004
$ws.Range("L135").Value = 13637.25
$ws.Range("M135").Value = -2469
$ws.Range("N135").Value = -18707.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 5339.3335
$ws.Range("I97").Value = 4405
$ws.Range("K97").Value = 4405
$ws.Range("M97").Value = -3909

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1779.9
$ws.Range("I61").Value = 1533.2222
$ws.Range("J61").Value = 4000
$ws.Range("K61").Value = 1533.2222
$ws.Range("L61").Value = 4000
$ws.Range("M61").Value = -1331.2222
$ws.Range("N61").Value = -4404

$ws.Range("H64").Value = 24000
$ws.Range("J64").Value = 24000
$ws.Range("L64").Value = 24000
$ws.Range("N64").Value = -24450

$ws.Range("H67").Value = 24000
$ws.Range("J67").Value = 24000
$ws.Range("L67").Value = 24000
$ws.Range("N67").Value = -25560

$ws.Range("H108").Value = 25000
$ws.Range("J108").Value = 25000
$ws.Range("L108").Value = 25000
$ws.Range("N108").Value = -32680

$ws.Range("H113").Value = 1779.9
$ws.Range("I113").Value = 1533.2222
$ws.Range("J113").Value = 4000
$ws.Range("K113").Value = 1533.2222
$ws.Range("L113").Value = 4000
$ws.Range("M113").Value = 636.7778000000001
$ws.Range("N113").Value = -8340

$ws.Range("H122").Value = 3275.2354
$ws.Range("I122").Value = 2310.3333
$ws.Range("J122").Value = 4360.75
$ws.Range("K122").Value = 6930.999899999999
$ws.Range("L122").Value = 13082.25
$ws.Range("M122").Value = -4480.999899999999
$ws.Range("N122").Value = -17982.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 30000
$ws.Range("J63").Value = 30000
$ws.Range("L63").Value = 30000
$ws.Range("N63").Value = -31248

$ws.Range("H66").Value = 30000
$ws.Range("J66").Value = 30000
$ws.Range("L66").Value = 90000
$ws.Range("N66").Value = -96240
